# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-27 14:24:48
#
# Normalizes the ordering of names inside the "Recorded By" column (G) on
# the "Session Analysis Results" sheet: the comma-separated list of
# recorders is re-sorted so that "System" sorts ahead of the other
# (lowercase-starting) entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

# Map of exact old "Recorded By" strings to their replacement, matching the
# canonical reordering applied across the sheet.
$replacements = @{
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2
    if ($replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
